$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header columns, copying the style (formats only) of the existing header cell
$ws.Range("F1").Copy() | Out-Null
$ws.Range("G1:H1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("G1").Value = "Round_4"
$ws.Range("H1").Value = "Round_5"

# Update the Round_1..Round_3 values and set new Round_4/Round_5 values
$data = @(
    @(3, 10, 2, 5, 4),
    @(2, 4, 0, 2, 8),
    @(7, 3, 9, 8, 5),
    @(5, 8, 4, 5, 4),
    @(10, 10, 1, 4, 6)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $vals = $data[$i]
    $ws.Cells.Item($row, 4).Value = $vals[0]
    $ws.Cells.Item($row, 5).Value = $vals[1]
    $ws.Cells.Item($row, 6).Value = $vals[2]
    $ws.Cells.Item($row, 7).Value = $vals[3]
    $ws.Cells.Item($row, 8).Value = $vals[4]
}
